# Applies a cyclic re-shuffle of species-record data between rows 22-25
# (22 <- 23 <- 24 <- 25 <- 22) and a swap between rows 32 and 33 on the
# "Artfynd" worksheet, matching the upstream data refresh captured in the
# commit diff. Column letters map to indices: A=1 B=2 D=4 E=5 F=6 G=7 H=8
# Q=17 R=18 (C, I and other columns are unchanged between these rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($ws, $row) {
    return [PSCustomObject]@{
        A = $ws.Cells.Item($row, 1).Value2
        B = $ws.Cells.Item($row, 2).Value2
        D = $ws.Cells.Item($row, 4).Value2
        E = $ws.Cells.Item($row, 5).Value2
        F = $ws.Cells.Item($row, 6).Value2
        G = $ws.Cells.Item($row, 7).Value2
        H = $ws.Cells.Item($row, 8).Value2
        Q = $ws.Cells.Item($row, 17).Value2
        R = $ws.Cells.Item($row, 18).Value2
    }
}

function Set-RowData($ws, $row, $data) {
    $ws.Cells.Item($row, 1).Value2 = $data.A
    $ws.Cells.Item($row, 2).Value2 = $data.B
    $ws.Cells.Item($row, 4).Value2 = $data.D
    $ws.Cells.Item($row, 5).Value2 = $data.E
    $ws.Cells.Item($row, 6).Value2 = $data.F
    $ws.Cells.Item($row, 7).Value2 = $data.G
    $ws.Cells.Item($row, 8).Value2 = $data.H
    $ws.Cells.Item($row, 17).Value2 = $data.Q
    $ws.Cells.Item($row, 18).Value2 = $data.R
}

# Capture original values for rows 22-25 and 32-33 before overwriting.
$row22 = Get-RowData $ws 22
$row23 = Get-RowData $ws 23
$row24 = Get-RowData $ws 24
$row25 = Get-RowData $ws 25
$row32 = Get-RowData $ws 32
$row33 = Get-RowData $ws 33

# Cyclic shift: 22 <- 23 <- 24 <- 25 <- 22(original)
Set-RowData $ws 22 $row23
Set-RowData $ws 23 $row24
Set-RowData $ws 24 $row25
Set-RowData $ws 25 $row22

# Swap rows 32 and 33
Set-RowData $ws 32 $row33
Set-RowData $ws 33 $row32
